# 2.1.1.1e - add 2023 (column Q) data, adjust row heights, clear stale selection,
# and turn off iterative calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column Q: copy formatting from the corresponding column P cell (so no
#    new cell styles are created), then write the 2023 values.
# ---------------------------------------------------------------------------

$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2023

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 74.605426356589135

$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 118.8

$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = 71.61643835616438

$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 95.703125

$ws.Range("P9").Copy()
$ws.Range("Q9").PasteSpecial(-4122)
$ws.Range("Q9").Value = 113.91018619934282

$ws.Range("P10").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$ws.Range("Q10").Value = 108.21501014198785

$ws.Range("P11").Copy()
$ws.Range("Q11").PasteSpecial(-4122)
$ws.Range("Q11").Value = 165.26684164479443

$ws.Range("P12").Copy()
$ws.Range("Q12").PasteSpecial(-4122)
$ws.Range("Q12").Value = 48.504446240905416

$ws.Range("P13").Copy()
$ws.Range("Q13").PasteSpecial(-4122)
$ws.Range("Q13").Value = 97.361348644026393

$ws.Range("P14").Copy()
$ws.Range("Q14").PasteSpecial(-4122)
$ws.Range("Q14").Value = 52.747252747252752

# ---------------------------------------------------------------------------
# 2. Row height adjustments.
# ---------------------------------------------------------------------------

$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 16.5
$ws.Rows.Item(8).RowHeight = 16.5
$ws.Rows.Item(9).RowHeight = 16.5
$ws.Rows.Item(10).RowHeight = 16.5
$ws.Rows.Item(11).RowHeight = 16.5
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(13).RowHeight = 16.5
$ws.Rows.Item(14).RowHeight = 16.5

# ---------------------------------------------------------------------------
# 3. Reset the lingering selection (was parked on R1, outside the used range)
#    back onto the sheet's data so it no longer points past column Q.
# ---------------------------------------------------------------------------

$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# 4. Turn off iterative calculation (workbook previously had manual iterate
#    settings enabled).
# ---------------------------------------------------------------------------

$excel.Iterate = $false
